# Conversion Script / Set_filter_file.xlsx - debugging fixes
# - Region_selection: clear the "1" flags that were set on rows 8-37 (set to 0)
# - Region_selection: append 5 new regions (DolAmroth, Gondor, Rohan, Harad, Mordor)
#   as selected (1) rows 56-60
# - Make Region_selection the active/selected sheet (was Year_selection)

$wb = $excel.ActiveWorkbook

$regionWs = $wb.Worksheets.Item("Region_selection")

# Un-flag rows 8 through 37 (values were 1, now 0)
$regionWs.Range("B8:B37").Value = 0

# Append the five new region rows
$regionWs.Range("A56").Value = "DolAmroth"
$regionWs.Range("A57").Value = "Gondor"
$regionWs.Range("A58").Value = "Rohan"
$regionWs.Range("A59").Value = "Harad"
$regionWs.Range("A60").Value = "Mordor"
$regionWs.Range("B56:B60").Value = 1

# Switch the active sheet/selection from Year_selection back to Region_selection
$regionWs.Activate()
$regionWs.Range("B61").Select()
